$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
  $c = $ws.Cells.Item($row, $col)
  $c.Value2 = "'" + $text
  $c.Style = "Normal"
}

# Update T2 (Stueckzinsen) value
Set-TextCell 2 20 "3.736986%"

# Use row 2 as a template: copy its values (no formatting, no number
# re-interpretation, no auto row-height recalculation) into the three
# new rows, then overwrite every column except Q (identical boilerplate
# text already correct after the copy) with the real data.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).PasteSpecial(-4163)
$ws.Rows.Item(4).PasteSpecial(-4163)
$ws.Rows.Item(5).PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Row 3 - Schweizerische Eidgenossenschaft
Set-TextCell 3 1  "Schweizerische Eidgenossenschaft"
Set-TextCell 3 2  "0.5 EIDG 30"
Set-TextCell 3 3  "CH"
Set-TextCell 3 4  "27.05.2015"
Set-TextCell 3 5  "27.05.2030"
Set-TextCell 3 6  "315.377 Mio"
Set-TextCell 3 7  "103.00"
Set-TextCell 3 8  "Nominal"
Set-TextCell 3 9  "1'000"
Set-TextCell 3 10 "100.00"
Set-TextCell 3 11 "CHF"
Set-TextCell 3 12 "CH0224397171"
Set-TextCell 3 13 "22439717"
Set-TextCell 3 14 "Nein"
Set-TextCell 3 15 "Nein"
Set-TextCell 3 16 "0.50%"
Set-TextCell 3 18 "26.05."
Set-TextCell 3 19 "30/360"
Set-TextCell 3 20 "0.338889%"

# Row 4 - Pfandbriefbank Schweizerischer Hypothekarinstitute AG
Set-TextCell 4 1  "Pfandbriefbank Schweizerischer Hypothekarinstitute AG"
Set-TextCell 4 2  "1.125 PB 23 S576"
Set-TextCell 4 3  "CH"
Set-TextCell 4 4  "16.11.2012"
Set-TextCell 4 5  "16.11.2023"
Set-TextCell 4 6  "130 Mio"
Set-TextCell 4 7  "101.146"
Set-TextCell 4 8  "Nominal"
Set-TextCell 4 9  "5'000"
Set-TextCell 4 10 "100.00"
Set-TextCell 4 11 "CHF"
Set-TextCell 4 12 "CH0199589588"
Set-TextCell 4 13 "19958958"
Set-TextCell 4 14 "Nein"
Set-TextCell 4 15 "Nein"
Set-TextCell 4 16 "1.125%"
Set-TextCell 4 18 "15.11."
Set-TextCell 4 19 "30/360"
Set-TextCell 4 20 "0.234375%"

# Row 5 - Pfandbriefzentrale der schweizerischen Kantonalbanken AG
Set-TextCell 5 1  "Pfandbriefzentrale der schweizerischen Kantonalbanken AG"
Set-TextCell 5 2  "1 PZ 23 S427"
Set-TextCell 5 3  "CH"
Set-TextCell 5 4  "13.11.2012"
Set-TextCell 5 5  "13.02.2023"
Set-TextCell 5 6  "295 Mio"
Set-TextCell 5 7  "100.934"
Set-TextCell 5 8  "Nominal"
Set-TextCell 5 9  "5'000"
Set-TextCell 5 10 "100.00"
Set-TextCell 5 11 "CHF"
Set-TextCell 5 12 "CH0198800325"
Set-TextCell 5 13 "19880032"
Set-TextCell 5 14 "Nein"
Set-TextCell 5 15 "Nein"
Set-TextCell 5 16 "1.00%"
Set-TextCell 5 18 "12.02."
Set-TextCell 5 19 "30/360"
Set-TextCell 5 20 "0.966667%"

Write-Host "Applied refdata updates: T2 changed, rows 3-5 added."
